$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 328.18182
$ws.Range("I41").Value = 222.5
$ws.Range("J41").Value = 455
$ws.Range("K41").Value = 222.5
$ws.Range("L41").Value = 455
$ws.Range("M41").Value = 217.5
$ws.Range("N41").Value = -1335
$ws.Range("H62").Value = 2735
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 2735
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H132").Value = 1048.375
$ws.Range("I132").Value = 998.3333
$ws.Range("K132").Value = 2994.9999
$ws.Range("M132").Value = -464.9998999999998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3756.4644
$ws.Range("I32").Value = 2321.558
$ws.Range("J32").Value = 8502.691999999999
$ws.Range("K32").Value = 2321.558
$ws.Range("L32").Value = 8502.691999999999
$ws.Range("M32").Value = -2034.558
$ws.Range("N32").Value = -9076.691999999999
$ws.Range("H122").Value = 1481.75
$ws.Range("I122").Value = 1434.4517
$ws.Range("K122").Value = 4303.355100000001
$ws.Range("M122").Value = -1853.355100000001
$ws.Range("H132").Value = 1852.7894
$ws.Range("I132").Value = 1482.6471
$ws.Range("K132").Value = 4447.9413
$ws.Range("M132").Value = -1917.9413
$ws.Range("H139").Value = 51999.5
$ws.Range("J139").Value = 51999.5
$ws.Range("L139").Value = 51999.5
$ws.Range("N139").Value = -62279.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8210.814
$ws.Range("I134").Value = 8583.862999999999
$ws.Range("K134").Value = 25751.589
$ws.Range("M134").Value = -23216.589
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3581.3635
$ws.Range("I31").Value = 1379.2727
$ws.Range("K31").Value = 1379.2727
$ws.Range("M31").Value = -1084.2727
$ws.Range("H34").Value = 3581.3635
$ws.Range("I34").Value = 1379.2727
$ws.Range("K34").Value = 1379.2727
$ws.Range("M34").Value = -1177.2727
$ws.Range("H50").Value = 17860
$ws.Range("J50").Value = 17860
$ws.Range("L50").Value = 17860
$ws.Range("N50").Value = -19110
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -37290
$ws.Range("H68").Value = 45000
$ws.Range("J68").Value = 45000
$ws.Range("L68").Value = 45000
$ws.Range("N68").Value = -46498
$ws.Range("H71").Value = 45000
$ws.Range("J71").Value = 45000
$ws.Range("L71").Value = 135000
$ws.Range("N71").Value = -142488
$ws.Range("H74").Value = 29999.666
$ws.Range("J74").Value = 29999.666
$ws.Range("L74").Value = 29999.666
$ws.Range("N74").Value = -31747.666
$ws.Range("H77").Value = 29999.666
$ws.Range("J77").Value = 29999.666
$ws.Range("L77").Value = 89998.99800000001
$ws.Range("N77").Value = -98734.99800000001
$ws.Range("H107").Value = 525.5789
$ws.Range("I107").Value = 436.5
$ws.Range("K107").Value = 436.5
$ws.Range("M107").Value = 1483.5
$ws.Range("H132").Value = 2706.3684
$ws.Range("I132").Value = 1385.3572
$ws.Range("K132").Value = 4156.071599999999
$ws.Range("M132").Value = -1626.071599999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 10395.25
$ws.Range("I56").Value = 10395.25
$ws.Range("K56").Value = 10395.25
$ws.Range("M56").Value = -9865.25
$ws.Range("H126").Value = 166669170
$ws.Range("I126").Value = 250000750
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 750002250
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -749997310
$ws.Range("N126").Value = -27880
$ws.Range("H131").Value = 9820300
$ws.Range("J131").Value = 19324.418
$ws.Range("L131").Value = 57973.254
$ws.Range("N131").Value = -68053.254
$ws.Range("H132").Value = 1574.0667
$ws.Range("J132").Value = 1602.8
$ws.Range("L132").Value = 14425.2
$ws.Range("N132").Value = -19485.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2393.3333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2393.3333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2393.3333
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -4389.3333
$ws.Range("H83").Value = 2393.3333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2393.3333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 11966.6665
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -21950.6665
$ws.Range("H113").Value = 1339.5834
$ws.Range("I113").Value = 1100.8334
$ws.Range("J113").Value = 1578.3334
$ws.Range("K113").Value = 1100.8334
$ws.Range("L113").Value = 1578.3334
$ws.Range("M113").Value = 1069.1666
$ws.Range("N113").Value = -5918.3334
$ws.Range("H122").Value = 1770.6207
$ws.Range("I122").Value = 1618.1052
$ws.Range("K122").Value = 4854.3156
$ws.Range("M122").Value = -2404.3156
$ws.Range("H126").Value = 55797.79
$ws.Range("I126").Value = 3385.5
$ws.Range("K126").Value = 10156.5
$ws.Range("M126").Value = -7686.5
$ws.Range("H132").Value = 3225.8965
$ws.Range("I132").Value = 2889.8333
$ws.Range("J132").Value = 4839
$ws.Range("K132").Value = 8669.499899999999
$ws.Range("L132").Value = 14517
$ws.Range("M132").Value = -6139.499899999999
$ws.Range("N132").Value = -19577
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 362.76923
$ws.Range("I55").Value = 129
$ws.Range("J55").Value = 563.1429000000001
$ws.Range("K55").Value = 129
$ws.Range("L55").Value = 563.1429000000001
$ws.Range("M55").Value = 44
$ws.Range("N55").Value = -909.1429000000001
$ws.Range("H61").Value = 3238.8
$ws.Range("I61").Value = 2798.5
$ws.Range("K61").Value = 2798.5
$ws.Range("M61").Value = -2596.5
$ws.Range("H113").Value = 3238.8
$ws.Range("I113").Value = 2798.5
$ws.Range("K113").Value = 2798.5
$ws.Range("M113").Value = -628.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 882.6667
$ws.Range("J113").Value = 959.2
$ws.Range("L113").Value = 2877.6
$ws.Range("N113").Value = -7217.6
$ws.Range("H132").Value = 3850.65
$ws.Range("I132").Value = 1319.6364
$ws.Range("J132").Value = 6944.1113
$ws.Range("K132").Value = 3958.9092
$ws.Range("L132").Value = 20832.3339
$ws.Range("M132").Value = -1428.9092
$ws.Range("N132").Value = -25892.3339
$ws.Range("H136").Value = 3613
$ws.Range("I136").Value = 3592.4167
$ws.Range("J136").Value = 3635.4546
$ws.Range("K136").Value = 10777.2501
$ws.Range("L136").Value = 10906.3638
$ws.Range("M136").Value = -8227.250100000001
$ws.Range("N136").Value = -16006.3638
